$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "FilesTab" query in B6 is rewritten: it now unions two sub-queries
# (sequencing files for Male participants + methylation array files for
# Female participants), adds a "File Category" column driven off the file
# type, renames a couple of source columns, and adds "Library Source
# Molecule". Rebuild the cell's text to match.
$filesTabQuery = @'
SELECT DISTINCT
    sqf.file_name AS "File Name",
    CASE
        WHEN sqf.type = 'sequencing_file' THEN 'Sequencing'
        ELSE sqf.type
    END AS "File Category",
    COALESCE(sqf.file_description, '') AS "File Description",
    sqf.file_type AS "File Type",
    CASE     
        WHEN sqf.file_size >= 1024 * 1024 * 1024 THEN 
            ROUND(sqf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB' 
        WHEN sqf.file_size >= 1024 * 1024 THEN 
            ROUND(sqf.file_size / (1024.0 * 1024.0), 2) || ' MB' 
        WHEN sqf.file_size >= 1024 THEN 
            ROUND(sqf.file_size / 1024.0, 2) || ' KB' 
        ELSE 
            ROUND(sqf.file_size, 2) || ' Bytes' 
    END AS "File Size",
    std.study AS "Study ID",
    prt.participant_id AS "Participant ID",
    smp.sample_id AS "Sample ID",
    sqf.dcf_indexd_guid AS "GUID",
    sqf.md5sum AS "MD5sum",
    COALESCE(sqf.library_selection, '') AS "Library Selection",
    COALESCE(sqf.library_source_material, '') AS "Library Source",
    COALESCE(sqf.library_strategy, '') AS "Library Strategy",
    COALESCE(sqf.library_source_molecule, '') AS "Library Source Molecule"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
JOIN 
    df_sequencing_file sqf ON smp.id = sqf."sample.id"
WHERE 
    std.study = 'phs002504' AND prt.sex_at_birth = 'Male'
UNION
SELECT DISTINCT
    maf.file_name AS "File Name",
    CASE
        WHEN maf.type = 'methylation_array_file' THEN 'Methylation array'
        ELSE maf.type
    END AS "File Category",
    COALESCE(maf.file_description, '') AS "File Description",
    maf.file_type AS "File Type",
    CASE     
        WHEN maf.file_size >= 1024 * 1024 * 1024 THEN 
            ROUND(maf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB' 
        WHEN maf.file_size >= 1024 * 1024 THEN 
            ROUND(maf.file_size / (1024.0 * 1024.0), 2) || ' MB' 
        WHEN maf.file_size >= 1024 THEN 
            ROUND(maf.file_size / 1024.0, 2) || ' KB' 
        ELSE 
            ROUND(maf.file_size, 2) || ' Bytes' 
    END AS "File Size",
    std.study AS "Study ID",
    prt.participant_id AS "Participant ID",
    smp.sample_id AS "Sample ID",
    maf.dcf_indexd_guid AS "GUID",
    maf.md5sum AS "MD5sum",
    '' AS "Library Selection",
    '' AS "Library Source",
    '' AS "Library Strategy",
    '' AS "Library Source Molecule"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_sample smp ON prt.id = smp."participant.id"
JOIN 
    df_methylation_array_file maf ON smp.id = maf."sample.id"
WHERE 
    std.study = 'phs002504' AND prt.sex_at_birth = 'Female'
ORDER BY 
    sqf.file_name ASC
LIMIT 100;
'@

$ws.Range("B6").Value = $filesTabQuery

# Setting that much text in a wrap-text cell can blow the row height far
# past Excel's normal autofit ceiling; keep row 6 at the same height the
# rest of the wrapped/overflowing query cells (rows 3-5) use.
$ws.Rows.Item(6).RowHeight = 409.5
